# Add a "Sprint 3" section to the timesheet, mirroring the existing
# "Sprint 2" / "Sprint 2 Notes" block structure (new rows 18-20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 18: "Sprint 3" section header (mirrors row 15 "Sprint 2") ---
$ws.Range("A18:B18").Merge()
$ws.Range("A15:B15").Copy()
$ws.Range("A18:B18").PasteSpecial(-4122)
$ws.Range("A18").Value = "Sprint 3"
$ws.Rows.Item(18).RowHeight = 39.75

# --- Row 19: blank data row under the Sprint 3 header (mirrors row 13,
#     the blank row with a top border / centered alignment that sits
#     below a sprint header) ---
$ws.Range("A19:B19").Merge()
$ws.Range("A13:B13").Copy()
$ws.Range("A19:B19").PasteSpecial(-4122)
$ws.Range("E19").Value = 6
$ws.Range("E19").HorizontalAlignment = -4108
$ws.Range("E19").VerticalAlignment = -4108

# --- Row 20: "Sprint 3 Notes" label + note content (mirrors row 17
#     "Sprint 2 Notes") ---
$ws.Range("A20:B20").Merge()
$ws.Range("A20:B20").HorizontalAlignment = -4108
$ws.Range("A20").Value = "Sprint 3 Notes"
$ws.Range("E20").Value = "Save and Load functionality, Testing"
$ws.Range("E20").HorizontalAlignment = -4108
$ws.Range("E20").VerticalAlignment = -4108

$ws.Range("E22").Select()
